$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 28.74177234654678
$ws.Range("C2").Value = 62.18064862231948
$ws.Range("D2").Value = 59.41903723182529
$ws.Range("E2").Value = 56.47178827322775
$ws.Range("F2").Value = 54.37067399358619
$ws.Range("G2").Value = 58.23525148940115
$ws.Range("H2").Value = 56.8477302638232
$ws.Range("I2").Value = 57.17064554575204
$ws.Range("J2").Value = 56.64290622476072
$ws.Range("K2").Value = 54.02713294622769
$ws.Range("L2").Value = 55.0292090577191
$ws.Range("M2").Value = 52.88930092093869
$ws.Range("N2").Value = 54.93833422353718
$ws.Range("O2").Value = 49.72986500935158
$ws.Range("B3").Value = 8.177571976521424
$ws.Range("C3").Value = 7.717323582603214
$ws.Range("D3").Value = 5.432532926446357
$ws.Range("E3").Value = 7.17264682455721
$ws.Range("F3").Value = 13.11085824781559
$ws.Range("G3").Value = 7.125826108612418
$ws.Range("H3").Value = 10.52177842172381
$ws.Range("I3").Value = 14.69025098970324
$ws.Range("J3").Value = 6.787406873057996
$ws.Range("K3").Value = 11.53477181475749
$ws.Range("L3").Value = 11.18889694508616
$ws.Range("M3").Value = 8.801851974533985
$ws.Range("N3").Value = 12.74357949037417
$ws.Range("O3").Value = 13.87163490493413
$ws.Range("B4").Value = 22.21252422189228
$ws.Range("C4").Value = 7.973034081930154
$ws.Range("D4").Value = 6.754282704207861
$ws.Range("E4").Value = 8.673540054460503
$ws.Range("F4").Value = 11.94596366796696
$ws.Range("G4").Value = 9.50825872596071
$ws.Range("H4").Value = 7.996585021822852
$ws.Range("I4").Value = 7.487498588276003
$ws.Range("J4").Value = 9.868069344961265
$ws.Range("K4").Value = 6.556537379690409
$ws.Range("L4").Value = 8.982363610595549
$ws.Range("M4").Value = 9.219360340659419
$ws.Range("N4").Value = 8.466025678127128
$ws.Range("O4").Value = 8.968023843363992
$ws.Range("B5").Value = 5.592130196801711
$ws.Range("C5").Value = 8.613168193627221
$ws.Range("D5").Value = 11.19887735205916
$ws.Range("E5").Value = 12.92390102608069
$ws.Range("F5").Value = 7.361311505706953
$ws.Range("G5").Value = 11.76215236634313
$ws.Range("H5").Value = 8.717483074968227
$ws.Range("I5").Value = 8.440789766146974
$ws.Range("J5").Value = 10.81433526362367
$ws.Range("K5").Value = 10.51400438927699
$ws.Range("L5").Value = 8.203068532160108
$ws.Range("M5").Value = 11.00924269940958
$ws.Range("N5").Value = 10.50127484650911
$ws.Range("O5").Value = 8.742318686486442
$ws.Range("B6").Value = 1.369700416411636
$ws.Range("C6").Value = 5.423932072986627
$ws.Range("D6").Value = 4.815798500377834
$ws.Range("E6").Value = 3.636701171525891
$ws.Range("F6").Value = 4.706933792974422
$ws.Range("G6").Value = 4.385944505072188
$ws.Range("H6").Value = 4.14850969698549
$ws.Range("I6").Value = 3.472439791167067
$ws.Range("J6").Value = 2.658535313986443
$ws.Range("K6").Value = 7.741315760596895
$ws.Range("L6").Value = 5.773136766474688
$ws.Range("M6").Value = 5.251133866213949
$ws.Range("N6").Value = 3.538084223655596
$ws.Range("O6").Value = 7.415161582107641
$ws.Range("B7").Value = 12.21067561036234
$ws.Range("C7").Value = 2.276477792435483
$ws.Range("D7").Value = 1.880934055234448
$ws.Range("E7").Value = 3.061119941322288
$ws.Range("F7").Value = 2.156512637645154
$ws.Range("G7").Value = 2.597370217502268
$ws.Range("H7").Value = 2.863384926906305
$ws.Range("I7").Value = 3.237132041241274
$ws.Range("J7").Value = 3.673140025492056
$ws.Range("K7").Value = 2.797756093461756
$ws.Range("L7").Value = 2.884885461452212
$ws.Range("M7").Value = 4.375475637945296
$ws.Range("N7").Value = 2.293580387922897
$ws.Range("O7").Value = 2.517121246257346
$ws.Range("B8").Value = 3.331198774662425
$ws.Range("C8").Value = 1.38509443490425
$ws.Range("D8").Value = 1.959204555567077
$ws.Range("E8").Value = 2.564092389800277
$ws.Range("F8").Value = 1.670571266212922
$ws.Range("G8").Value = 1.911375088283461
$ws.Range("H8").Value = 2.041588219446828
$ws.Range("I8").Value = 1.880049123742429
$ws.Range("J8").Value = 2.566124269033158
$ws.Range("K8").Value = 2.16288470023133
$ws.Range("L8").Value = 3.554673745957689
$ws.Range("M8").Value = 2.998704105299218
$ws.Range("N8").Value = 2.811792779360872
$ws.Range("O8").Value = 1.938378041965687
$ws.Range("B9").Value = 18.3644264568014
$ws.Range("C9").Value = 4.430321219193573
$ws.Range("D9").Value = 8.539332674281985
$ws.Range("E9").Value = 5.496210319025387
$ws.Range("F9").Value = 4.67717488809181
$ws.Range("G9").Value = 4.473821498824696
$ws.Range("H9").Value = 6.86294037432328
$ws.Range("I9").Value = 3.62119415397097
$ws.Range("J9").Value = 6.989482685084679
$ws.Range("K9").Value = 4.665596915757445
$ws.Range("L9").Value = 4.383765880554509
$ws.Range("M9").Value = 5.454930454999854
$ws.Range("N9").Value = 4.707328370513061
$ws.Range("O9").Value = 6.817496685533188
